$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "97.582.70"
$ws.Range("E2").Value = "  -1.69%  "
$ws.Range("D3").Value = "3.413.51"
$ws.Range("E3").Value = "  +3.90%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "254.75"
$ws.Range("E5").Value = "  +0.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "653.77"
$ws.Range("E6").Value = "  +5.02%  "
$ws.Range("E7").Value = "  +3.28%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.424"
$ws.Range("E8").Value = "  +6.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.07"
$ws.Range("E9").Value = "  +10.30%  "
$ws.Range("E10").Value = "  -0.06%  "
$ws.Range("D11").Value = "3.409.50"
$ws.Range("E11").Value = "  +3.88%  "
$ws.Range("E12").Value = "  +6.07%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "41.86"
$ws.Range("E13").Value = "  +5.96%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.28"
$ws.Range("E14").Value = "  +14.48%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000260"
$ws.Range("E15").Value = "  +5.17%  "
$ws.Range("D16").Value = "97.263.28"
$ws.Range("E16").Value = "  -1.85%  "
$ws.Range("D17").Value = "4.046.44"
$ws.Range("E17").Value = "  +4.46%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.62"
$ws.Range("E18").Value = "  +35.72%  "
$ws.Range("D19").Value = "3.413.28"
$ws.Range("E19").Value = "  +4.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.54"
$ws.Range("E20").Value = "  +14.99%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.505"
$ws.Range("E21").Value = "  +53.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.92"
$ws.Range("E22").Value = "  +17.63%  "
$ws.Range("E23").Value = "  +0.56%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "506.52"
$ws.Range("E24").Value = "  +4.01%  "
$ws.Range("E25").Value = "  +2.29%  "
$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "98.79"
$ws.Range("E26").Value = "  +10.81%  "
$ws.Range("B27").Value = "NEARProtocol"
$ws.Range("C27").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.14"
$ws.Range("E27").Value = "  +9.40%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.76"
$ws.Range("E28").Value = "  +6.44%  "
$ws.Range("D29").Value = "3.597.40"
$ws.Range("E29").Value = "  +5.04%  "
$ws.Range("E30").Value = "  +14.95%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.43"
$ws.Range("E31").Value = "  +10.26%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.199"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  -0.11%  "
$ws.Range("E34").Value = "  +0.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.575"
$ws.Range("E35").Value = "  +21.34%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "29.92"
$ws.Range("E36").Value = "  +6.97%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.27"
$ws.Range("E37").Value = "  +17.56%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.75"
$ws.Range("E38").Value = "  +7.74%  "
$ws.Range("E39").Value = "  +2.80%  "
$ws.Range("E40").Value = "  +15.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "513.16"
$ws.Range("E41").Value = "  +5.59%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "24.71"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.859"
$ws.Range("E43").Value = "  +11.35%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0422"
$ws.Range("E44").Value = "  +24.96%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.66"
$ws.Range("E45").Value = "  -1.27%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.51"
$ws.Range("E46").Value = "  +16.63%  "
$ws.Range("E47").Value = "  +6.17%  "
$ws.Range("E48").Value = "  +0.02%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.21"
$ws.Range("E49").Value = "  +13.18%  "
$ws.Range("E50").Value = "  +16.54%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.08"
$ws.Range("E51").Value = "  +7.34%  "
